$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the date-serial values in column A (rows 2-22) with quarterly
# period labels (e.g. "2004Q4") stored as plain text, and make them share
# the same (text) style as the header row, instead of the removed
# custom date-time number format.
$startYear = 2004
for ($row = 2; $row -le 22; $row++) {
    $year = $startYear + ($row - 2)
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "$($year)Q4"
    $cell.Style = $ws.Cells.Item(1, 1).Style
}
